$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-space the shift boundaries for rows 2-5 (3h -> 4h/3h/3h/3h blocks),
# so the last (5th) shift of the day (row 6, 20:00-21:00) is no longer needed.
$ws.Range("D2").Value = 45292.5
$ws.Range("C3").Value = 45292.5
$ws.Range("D3").Value = 45292.625
$ws.Range("C4").Value = 45292.625
$ws.Range("D4").Value = 45292.75
$ws.Range("C5").Value = 45292.75
$ws.Range("D5").Value = 45292.875

# Drop the now-unused 5th shift row entirely; this shifts every following
# (already-blank) row up by one, which also removes the trailing blank
# row 22 and shrinks the sheet's used range to A1:E21.
$ws.Rows(6).Delete()
